$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows at 1183 (pushes existing 1183:1205 down to 1189:1211),
# copying the formatting (incl. the date style on column D) from the row above.
$ws.Rows("1183:1188").Insert()

$ws.Range('A1183').Value = 6
$ws.Range('B1183').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C1183').Value = 'Metropolitana'
$ws.Range('D1183').Value = 44595
$ws.Range('E1183').Value = 13
$ws.Range('F1183').Value = 100114013
$ws.Range('G1183').Value = 'Zanahoria'
$ws.Range('H1183').Value = 'Sin especificar'
$ws.Range('I1183').Value = 'Camote'
$ws.Range('J1183').Value = 330
$ws.Range('K1183').Value = 5000
$ws.Range('L1183').Value = 5000
$ws.Range('M1183').Value = 5000
$ws.Range('N1183').Value = '$/saco 20 kilos'
$ws.Range('O1183').Value = 'Región Metropolitana'
$ws.Range('P1183').Value = 250
$ws.Range('Q1183').Value = 20
$ws.Range('R1183').Value = 'Hortaliza'

$ws.Range('A1184').Value = 6
$ws.Range('B1184').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C1184').Value = 'Metropolitana'
$ws.Range('D1184').Value = 44595
$ws.Range('E1184').Value = 13
$ws.Range('F1184').Value = 100114013
$ws.Range('G1184').Value = 'Zanahoria'
$ws.Range('H1184').Value = 'Sin especificar'
$ws.Range('I1184').Value = 'Camote'
$ws.Range('J1184').Value = 400
$ws.Range('K1184').Value = 5000
$ws.Range('L1184').Value = 5000
$ws.Range('M1184').Value = 5000
$ws.Range('N1184').Value = '$/saco 20 kilos'
$ws.Range('O1184').Value = 'Región de Ñuble'
$ws.Range('P1184').Value = 250
$ws.Range('Q1184').Value = 20
$ws.Range('R1184').Value = 'Hortaliza'

$ws.Range('A1185').Value = 6
$ws.Range('B1185').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C1185').Value = 'Metropolitana'
$ws.Range('D1185').Value = 44595
$ws.Range('E1185').Value = 13
$ws.Range('F1185').Value = 100114013
$ws.Range('G1185').Value = 'Zanahoria'
$ws.Range('H1185').Value = 'Sin especificar'
$ws.Range('I1185').Value = 'Primera'
$ws.Range('J1185').Value = 2300
$ws.Range('K1185').Value = 6000
$ws.Range('L1185').Value = 6500
$ws.Range('M1185').Value = 6196
$ws.Range('N1185').Value = '$/saco 20 kilos'
$ws.Range('O1185').Value = 'Región Metropolitana'
$ws.Range('P1185').Value = 310
$ws.Range('Q1185').Value = 20
$ws.Range('R1185').Value = 'Hortaliza'

$ws.Range('A1186').Value = 6
$ws.Range('B1186').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C1186').Value = 'Metropolitana'
$ws.Range('D1186').Value = 44595
$ws.Range('E1186').Value = 13
$ws.Range('F1186').Value = 100114013
$ws.Range('G1186').Value = 'Zanahoria'
$ws.Range('H1186').Value = 'Sin especificar'
$ws.Range('I1186').Value = 'Primera'
$ws.Range('J1186').Value = 2800
$ws.Range('K1186').Value = 6000
$ws.Range('L1186').Value = 6500
$ws.Range('M1186').Value = 6232
$ws.Range('N1186').Value = '$/saco 20 kilos'
$ws.Range('O1186').Value = 'Región de Ñuble'
$ws.Range('P1186').Value = 312
$ws.Range('Q1186').Value = 20
$ws.Range('R1186').Value = 'Hortaliza'

$ws.Range('A1187').Value = 6
$ws.Range('B1187').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C1187').Value = 'Metropolitana'
$ws.Range('D1187').Value = 44595
$ws.Range('E1187').Value = 13
$ws.Range('F1187').Value = 100114013
$ws.Range('G1187').Value = 'Zanahoria'
$ws.Range('H1187').Value = 'Sin especificar'
$ws.Range('I1187').Value = 'Segunda'
$ws.Range('J1187').Value = 500
$ws.Range('K1187').Value = 5500
$ws.Range('L1187').Value = 5500
$ws.Range('M1187').Value = 5500
$ws.Range('N1187').Value = '$/saco 20 kilos'
$ws.Range('O1187').Value = 'Región Metropolitana'
$ws.Range('P1187').Value = 275
$ws.Range('Q1187').Value = 20
$ws.Range('R1187').Value = 'Hortaliza'

$ws.Range('A1188').Value = 6
$ws.Range('B1188').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C1188').Value = 'Metropolitana'
$ws.Range('D1188').Value = 44595
$ws.Range('E1188').Value = 13
$ws.Range('F1188').Value = 100114013
$ws.Range('G1188').Value = 'Zanahoria'
$ws.Range('H1188').Value = 'Sin especificar'
$ws.Range('I1188').Value = 'Segunda'
$ws.Range('J1188').Value = 400
$ws.Range('K1188').Value = 5000
$ws.Range('L1188').Value = 5000
$ws.Range('M1188').Value = 5000
$ws.Range('N1188').Value = '$/saco 20 kilos'
$ws.Range('O1188').Value = 'Región de Ñuble'
$ws.Range('P1188').Value = 250
$ws.Range('Q1188').Value = 20
$ws.Range('R1188').Value = 'Hortaliza'

